$d = $word.ActiveDocument

# 1. "For a long time have models" -> "For a long time, models"
$d.Content.Find.Execute("For a long time have models", $true, $false, $false, $false, $false, $true, 1, $false, "For a long time, models", 2) | Out-Null

# 2. "including gravitationally and seismologically." -> "including gravitational and seismological studies."
$d.Content.Find.Execute("including gravitationally and seismologically.", $true, $false, $false, $false, $false, $true, 1, $false, "including gravitational and seismological studies.", 2) | Out-Null

# 3. "interpolated but there is no way" -> "interpolated; there is no way"
$d.Content.Find.Execute("the seismic data has to be interpolated but there is no way", $true, $false, $false, $false, $false, $true, 1, $false, "the seismic data has to be interpolated; there is no way", 2) | Out-Null

# 4. "To get around this problem" -> "In order to determine a method's accuracy"
$d.Content.Find.Execute("true Moho depth model. To get around this problem a method of cross-validation", $true, $false, $false, $false, $false, $true, 1, $false, "true Moho depth model. In order to determine a method’s accuracy a method of cross-validation", 2) | Out-Null

# 5. "will tell how good gravitational models are in regions" -> "will evaluate the accuracy of gravitational models in regions"
$d.Content.Find.Execute("The results from this cross-validation will tell how good gravitational models are in regions", $true, $false, $false, $false, $false, $true, 1, $false, "The results from this cross-validation will evaluate the accuracy of gravitational models in regions", 2) | Out-Null

# 6. Paraná Basin sentence rewrite
$d.Content.Find.Execute("Basin is thought to have large igneous intrusions resulting in a shallower Moho than expected, in an attempt to decrease the errors on the model this intrusion will be modelled. The cross-validation", $true, $false, $false, $false, $false, $true, 1, $false, "Basin, South America is thought to have large igneous intrusions resulting in a shallower Moho than expected. In this study these intrusions will be modelled in an attempt to decrease the errors on the Moho model currently used. The cross-validation", 2) | Out-Null

# 7. Give the "á" in "Paraná" its own run (mirrors the source formatting mark on that glyph)
$accentRange = $d.Content
$found = $accentRange.Find.Execute("á")
if ($found) {
    $accentRange.Font.Name = "Calibri"
}

Write-Output $d.Content.Text
